$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "70-11=59"
$t.Cell(1, 2).Range.Text = "48+27=75"
$t.Cell(1, 3).Range.Text = "48+14=62"
$t.Cell(1, 4).Range.Text = "66-28=38"
$t.Cell(1, 5).Range.Text = "70-5=65"
$t.Cell(2, 1).Range.Text = "83-55=28"
$t.Cell(2, 2).Range.Text = "74+18=92"
$t.Cell(2, 3).Range.Text = "56-49=7"
$t.Cell(2, 4).Range.Text = "81-66=15"
$t.Cell(2, 5).Range.Text = "64-48=16"
$t.Cell(3, 1).Range.Text = "9+66=75"
$t.Cell(3, 2).Range.Text = "7+75=82"
$t.Cell(3, 3).Range.Text = "46-37=9"
$t.Cell(3, 4).Range.Text = "41-39=2"
$t.Cell(3, 5).Range.Text = "73-15=58"
$t.Cell(4, 1).Range.Text = "94-65=29"
$t.Cell(4, 2).Range.Text = "61-13=48"
$t.Cell(4, 3).Range.Text = "13+78=91"
$t.Cell(4, 4).Range.Text = "31-13=18"
$t.Cell(4, 5).Range.Text = "91-39=52"
$t.Cell(5, 1).Range.Text = "87+5=92"
$t.Cell(5, 2).Range.Text = "60-58=2"
$t.Cell(5, 3).Range.Text = "79+14=93"
$t.Cell(5, 4).Range.Text = "8+27=35"
$t.Cell(5, 5).Range.Text = "68+15=83"
$t.Cell(6, 1).Range.Text = "65-17=48"
$t.Cell(6, 2).Range.Text = "18+43=61"
$t.Cell(6, 3).Range.Text = "7+46=53"
$t.Cell(6, 4).Range.Text = "19+22=41"
$t.Cell(6, 5).Range.Text = "5+88=93"
$t.Cell(7, 1).Range.Text = "29+12=41"
$t.Cell(7, 2).Range.Text = "82-68=14"
$t.Cell(7, 3).Range.Text = "80-21=59"
$t.Cell(7, 4).Range.Text = "66+15=81"
$t.Cell(7, 5).Range.Text = "90-77=13"
$t.Cell(8, 1).Range.Text = "93-24=69"
$t.Cell(8, 2).Range.Text = "85-68=17"
$t.Cell(8, 3).Range.Text = "27+37=64"
$t.Cell(8, 4).Range.Text = "65+27=92"
$t.Cell(8, 5).Range.Text = "19+49=68"
$t.Cell(9, 1).Range.Text = "17+14=31"
$t.Cell(9, 2).Range.Text = "9+49=58"
$t.Cell(9, 3).Range.Text = "95-76=19"
$t.Cell(9, 4).Range.Text = "8+79=87"
$t.Cell(9, 5).Range.Text = "28+69=97"
$t.Cell(10, 1).Range.Text = "47+45=92"
$t.Cell(10, 2).Range.Text = "41-19=22"
$t.Cell(10, 3).Range.Text = "30-13=17"
$t.Cell(10, 4).Range.Text = "20-7=13"
$t.Cell(10, 5).Range.Text = "76+9=85"
$t.Cell(11, 1).Range.Text = "59+39=98"
$t.Cell(11, 2).Range.Text = "49+33=82"
$t.Cell(11, 3).Range.Text = "17+45=62"
$t.Cell(11, 4).Range.Text = "69+19=88"
$t.Cell(11, 5).Range.Text = "28+36=64"
$t.Cell(12, 1).Range.Text = "68-49=19"
$t.Cell(12, 2).Range.Text = "33+19=52"
$t.Cell(12, 3).Range.Text = "48+34=82"
$t.Cell(12, 4).Range.Text = "64-37=27"
$t.Cell(12, 5).Range.Text = "97-48=49"
$t.Cell(13, 1).Range.Text = "56+16=72"
$t.Cell(13, 2).Range.Text = "64-6=58"
$t.Cell(13, 3).Range.Text = "46+9=55"
$t.Cell(13, 4).Range.Text = "91-35=56"
$t.Cell(13, 5).Range.Text = "37+56=93"
$t.Cell(14, 1).Range.Text = "86+6=92"
$t.Cell(14, 2).Range.Text = "82-13=69"
$t.Cell(14, 3).Range.Text = "96-39=57"
$t.Cell(14, 4).Range.Text = "42-39=3"
$t.Cell(14, 5).Range.Text = "54-29=25"
$t.Cell(15, 1).Range.Text = "8+43=51"
$t.Cell(15, 2).Range.Text = "16-9=7"
$t.Cell(15, 3).Range.Text = "92-73=19"
$t.Cell(15, 4).Range.Text = "88+5=93"
$t.Cell(15, 5).Range.Text = "75-16=59"
$t.Cell(16, 1).Range.Text = "24-17=7"
$t.Cell(16, 2).Range.Text = "75-38=37"
$t.Cell(16, 3).Range.Text = "9+77=86"
$t.Cell(16, 4).Range.Text = "84-75=9"
$t.Cell(16, 5).Range.Text = "84-48=36"
$t.Cell(17, 1).Range.Text = "50-1=49"
$t.Cell(17, 2).Range.Text = "26+49=75"
$t.Cell(17, 3).Range.Text = "12+19=31"
$t.Cell(17, 4).Range.Text = "44-25=19"
$t.Cell(17, 5).Range.Text = "49+18=67"
$t.Cell(18, 1).Range.Text = "49+3=52"
$t.Cell(18, 2).Range.Text = "81-73=8"
$t.Cell(18, 3).Range.Text = "83-17=66"
$t.Cell(18, 4).Range.Text = "17+58=75"
$t.Cell(18, 5).Range.Text = "82-43=39"
$t.Cell(19, 1).Range.Text = "58+26=84"
$t.Cell(19, 2).Range.Text = "56+7=63"
$t.Cell(19, 3).Range.Text = "66+29=95"
$t.Cell(19, 4).Range.Text = "92-54=38"
$t.Cell(19, 5).Range.Text = "42-9=33"
$t.Cell(20, 1).Range.Text = "95-39=56"
$t.Cell(20, 2).Range.Text = "19+54=73"
$t.Cell(20, 3).Range.Text = "5+68=73"
$t.Cell(20, 4).Range.Text = "96-69=27"
$t.Cell(20, 5).Range.Text = "38+34=72"
